$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.617202
$ws.Range("H2").Value = 82.851606
$ws.Range("I2").Value = 0.4561121035646509
$ws.Range("J2").Value = 0.4561121035646509
$ws.Range("M2").Value = 127.5808283333333
$ws.Range("N2").Value = 382.742485
$ws.Range("O2").Value = 0.1532286472569342
$ws.Range("P2").Value = 0.1532286472569342
$ws.Range("Q2").Value = 3523.42550740899
$ws.Range("R2").Value = 31710.82956668091
$ws.Range("S2").Value = 0.06988944062672614
$ws.Range("T2").Value = 0.06988944062672614

$ws.Range("G3").Value = 27.617202
$ws.Range("H3").Value = 82.851606
$ws.Range("I3").Value = 0.4561121035646509
$ws.Range("J3").Value = 0.4561121035646509
$ws.Range("O3").Value = 0.341528235684153
$ws.Range("P3").Value = 0.341528235684153
$ws.Range("Q3").Value = 7853.291918006395
$ws.Range("R3").Value = 70679.62726205755
$ws.Range("S3").Value = 0.1557751620046229
$ws.Range("T3").Value = 0.1557751620046229

$ws.Range("G4").Value = 27.617202
$ws.Range("H4").Value = 82.851606
$ws.Range("I4").Value = 0.4561121035646509
$ws.Range("J4").Value = 0.4561121035646509
$ws.Range("M4").Value = 155.9545746666667
$ws.Range("N4").Value = 467.863724
$ws.Range("O4").Value = 0.1873064223040503
$ws.Range("P4").Value = 0.1873064223040503
$ws.Range("Q4").Value = 4307.028991393417
$ws.Range("R4").Value = 38763.26092254075
$ws.Range("S4").Value = 0.08543272628826921
$ws.Range("T4").Value = 0.08543272628826923

$ws.Range("G5").Value = 27.617202
$ws.Range("H5").Value = 82.851606
$ws.Range("I5").Value = 0.4561121035646509
$ws.Range("J5").Value = 0.4561121035646509
$ws.Range("M5").Value = 132.7811556666666
$ws.Range("N5").Value = 398.343467
$ws.Range("O5").Value = 0.1594744063806953
$ws.Range("P5").Value = 0.1594744063806954
$ws.Range("Q5").Value = 3667.043997839778
$ws.Range("R5").Value = 33003.395980558
$ws.Range("S5").Value = 0.07273820695902294
$ws.Range("T5").Value = 0.07273820695902296

$ws.Range("G6").Value = 27.617202
$ws.Range("H6").Value = 82.851606
$ws.Range("I6").Value = 0.4561121035646509
$ws.Range("J6").Value = 0.4561121035646509
$ws.Range("M6").Value = 131.9384486666667
$ws.Range("N6").Value = 395.815346
$ws.Range("O6").Value = 0.1584622883741672
$ws.Range("P6").Value = 0.1584622883741672
$ws.Range("Q6").Value = 3643.770788393964
$ws.Range("R6").Value = 32793.93709554568
$ws.Range("S6").Value = 0.07227656768600972
$ws.Range("T6").Value = 0.07227656768600974

$ws.Range("I7").Value = 0.4154106045258149
$ws.Range("J7").Value = 0.4154106045258149
$ws.Range("M7").Value = 127.5808283333333
$ws.Range("N7").Value = 382.742485
$ws.Range("O7").Value = 0.1532286472569342
$ws.Range("P7").Value = 0.1532286472569342
$ws.Range("Q7").Value = 3209.010040723419
$ws.Range("R7").Value = 28881.09036651077
$ws.Range("S7").Value = 0.06365280498767588
$ws.Range("T7").Value = 0.06365280498767588

$ws.Range("I8").Value = 0.4154106045258149
$ws.Range("J8").Value = 0.4154106045258149
$ws.Range("O8").Value = 0.341528235684153
$ws.Range("P8").Value = 0.341528235684153
$ws.Range("S8").Value = 0.1418744508481889
$ws.Range("T8").Value = 0.141874450848189

$ws.Range("I9").Value = 0.4154106045258149
$ws.Range("J9").Value = 0.4154106045258149
$ws.Range("M9").Value = 155.9545746666667
$ws.Range("N9").Value = 467.863724
$ws.Range("O9").Value = 0.1873064223040503
$ws.Range("P9").Value = 0.1873064223040503
$ws.Range("Q9").Value = 3922.688091462463
$ws.Range("R9").Value = 35304.19282316217
$ws.Range("S9").Value = 0.07780907412089309
$ws.Range("T9").Value = 0.07780907412089311

$ws.Range("I10").Value = 0.4154106045258149
$ws.Range("J10").Value = 0.4154106045258149
$ws.Range("M10").Value = 132.7811556666666
$ws.Range("N10").Value = 398.343467
$ws.Range("O10").Value = 0.1594744063806953
$ws.Range("P10").Value = 0.1594744063806954
$ws.Range("Q10").Value = 3339.812629527077
$ws.Range("R10").Value = 30058.31366574369
$ws.Range("S10").Value = 0.06624735956100013
$ws.Range("T10").Value = 0.06624735956100014

$ws.Range("I11").Value = 0.4154106045258149
$ws.Range("J11").Value = 0.4154106045258149
$ws.Range("M11").Value = 131.9384486666667
$ws.Range("N11").Value = 395.815346
$ws.Range("O11").Value = 0.1584622883741672
$ws.Range("P11").Value = 0.1584622883741672
$ws.Range("Q11").Value = 3318.616222043952
$ws.Range("R11").Value = 29867.54599839557
$ws.Range("S11").Value = 0.06582691500805679
$ws.Range("T11").Value = 0.0658269150080568

$ws.Range("G12").Value = 0.3390106666666666
$ws.Range("H12").Value = 1.017032
$ws.Range("I12").Value = 0.005598933168689138
$ws.Range("J12").Value = 0.005598933168689139
$ws.Range("M12").Value = 127.5808283333333
$ws.Range("N12").Value = 382.742485
$ws.Range("O12").Value = 0.1532286472569342
$ws.Range("P12").Value = 0.1532286472569342
$ws.Range("Q12").Value = 43.25126166716888
$ws.Range("R12").Value = 389.2613550045199
$ws.Range("S12").Value = 0.0008579169555202167
$ws.Range("T12").Value = 0.0008579169555202168

$ws.Range("G13").Value = 0.3390106666666666
$ws.Range("H13").Value = 1.017032
$ws.Range("I13").Value = 0.005598933168689138
$ws.Range("J13").Value = 0.005598933168689139
$ws.Range("O13").Value = 0.341528235684153
$ws.Range("P13").Value = 0.341528235684153
$ws.Range("Q13").Value = 96.40186318119021
$ws.Range("R13").Value = 867.6167686307118
$ws.Range("S13").Value = 0.001912193766815885
$ws.Range("T13").Value = 0.001912193766815885

$ws.Range("G14").Value = 0.3390106666666666
$ws.Range("H14").Value = 1.017032
$ws.Range("I14").Value = 0.005598933168689138
$ws.Range("J14").Value = 0.005598933168689139
$ws.Range("M14").Value = 155.9545746666667
$ws.Range("N14").Value = 467.863724
$ws.Range("O14").Value = 0.1873064223040503
$ws.Range("P14").Value = 0.1873064223040503
$ws.Range("Q14").Value = 52.8702643274631
$ws.Range("R14").Value = 475.832378947168
$ws.Range("S14").Value = 0.001048716140546642
$ws.Range("T14").Value = 0.001048716140546642

$ws.Range("G15").Value = 0.3390106666666666
$ws.Range("H15").Value = 1.017032
$ws.Range("I15").Value = 0.005598933168689138
$ws.Range("J15").Value = 0.005598933168689139
$ws.Range("M15").Value = 132.7811556666666
$ws.Range("N15").Value = 398.343467
$ws.Range("O15").Value = 0.1594744063806953
$ws.Range("P15").Value = 0.1594744063806954
$ws.Range("Q15").Value = 45.0142281033271
$ws.Range("R15").Value = 405.128052929944
$ws.Range("S15").Value = 0.0008928865434418859
$ws.Range("T15").Value = 0.0008928865434418862

$ws.Range("G16").Value = 0.3390106666666666
$ws.Range("H16").Value = 1.017032
$ws.Range("I16").Value = 0.005598933168689138
$ws.Range("J16").Value = 0.005598933168689139
$ws.Range("M16").Value = 131.9384486666667
$ws.Range("N16").Value = 395.815346
$ws.Range("O16").Value = 0.1584622883741672
$ws.Range("P16").Value = 0.1584622883741672
$ws.Range("Q16").Value = 44.72854144145244
$ws.Range("R16").Value = 402.556872973072
$ws.Range("S16").Value = 0.0008872197623645077
$ws.Range("T16").Value = 0.0008872197623645081

$ws.Range("G17").Value = 6.958474666666667
$ws.Range("H17").Value = 20.875424
$ws.Range("I17").Value = 0.1149227397407842
$ws.Range("J17").Value = 0.1149227397407842
$ws.Range("M17").Value = 127.5808283333333
$ws.Range("N17").Value = 382.742485
$ws.Range("O17").Value = 0.1532286472569342
$ws.Range("P17").Value = 0.1532286472569342
$ws.Range("Q17").Value = 887.7679619098489
$ws.Range("R17").Value = 7989.91165718864
$ws.Range("S17").Value = 0.01760945594954108
$ws.Range("T17").Value = 0.01760945594954108

$ws.Range("G18").Value = 6.958474666666667
$ws.Range("H18").Value = 20.875424
$ws.Range("I18").Value = 0.1149227397407842
$ws.Range("J18").Value = 0.1149227397407842
$ws.Range("O18").Value = 0.341528235684153
$ws.Range("P18").Value = 0.341528235684153
$ws.Range("Q18").Value = 1978.728071778798
$ws.Range("R18").Value = 17808.55264600918
$ws.Range("S18").Value = 0.03924936054365914
$ws.Range("T18").Value = 0.03924936054365914

$ws.Range("G19").Value = 6.958474666666667
$ws.Range("H19").Value = 20.875424
$ws.Range("I19").Value = 0.1149227397407842
$ws.Range("J19").Value = 0.1149227397407842
$ws.Range("M19").Value = 155.9545746666667
$ws.Range("N19").Value = 467.863724
$ws.Range("O19").Value = 0.1873064223040503
$ws.Range("P19").Value = 0.1873064223040503
$ws.Range("Q19").Value = 1085.205956968775
$ws.Range("R19").Value = 9766.853612718976
$ws.Range("S19").Value = 0.02152576722222579
$ws.Range("T19").Value = 0.0215257672222258

$ws.Range("G20").Value = 6.958474666666667
$ws.Range("H20").Value = 20.875424
$ws.Range("I20").Value = 0.1149227397407842
$ws.Range("J20").Value = 0.1149227397407842
$ws.Range("M20").Value = 132.7811556666666
$ws.Range("N20").Value = 398.343467
$ws.Range("O20").Value = 0.1594744063806953
$ws.Range("P20").Value = 0.1594744063806954
$ws.Range("Q20").Value = 923.954307917223
$ws.Range("R20").Value = 8315.588771255007
$ws.Range("S20").Value = 0.01832723569980471
$ws.Range("T20").Value = 0.01832723569980472

$ws.Range("G21").Value = 6.958474666666667
$ws.Range("H21").Value = 20.875424
$ws.Range("I21").Value = 0.1149227397407842
$ws.Range("J21").Value = 0.1149227397407842
$ws.Range("M21").Value = 131.9384486666667
$ws.Range("N21").Value = 395.815346
$ws.Range("O21").Value = 0.1584622883741672
$ws.Range("P21").Value = 0.1584622883741672
$ws.Range("Q21").Value = 918.0903526063004
$ws.Range("R21").Value = 8262.813173456703
$ws.Range("S21").Value = 0.01821092032555352
$ws.Range("T21").Value = 0.01821092032555352

$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.481706
$ws.Range("H22").Value = 1.445118
$ws.Range("I22").Value = 0.007955619000060676
$ws.Range("J22").Value = 0.007955619000060676
$ws.Range("M22").Value = 127.5808283333333
$ws.Range("N22").Value = 382.742485
$ws.Range("O22").Value = 0.1532286472569342
$ws.Range("P22").Value = 0.1532286472569342
$ws.Range("Q22").Value = 61.45645049313666
$ws.Range("R22").Value = 553.1080544382299
$ws.Range("S22").Value = 0.001219028737470861
$ws.Range("T22").Value = 0.001219028737470861

$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0.481706
$ws.Range("H23").Value = 1.445118
$ws.Range("I23").Value = 0.007955619000060676
$ws.Range("J23").Value = 0.007955619000060676
$ws.Range("O23").Value = 0.341528235684153
$ws.Range("P23").Value = 0.341528235684153
$ws.Range("Q23").Value = 136.9790406955487
$ws.Range("R23").Value = 1232.811366259938
$ws.Range("S23").Value = 0.002717068520866048
$ws.Range("T23").Value = 0.002717068520866048

$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.481706
$ws.Range("H24").Value = 1.445118
$ws.Range("I24").Value = 0.007955619000060676
$ws.Range("J24").Value = 0.007955619000060676
$ws.Range("M24").Value = 155.9545746666667
$ws.Range("N24").Value = 467.863724
$ws.Range("O24").Value = 0.1873064223040503
$ws.Range("P24").Value = 0.1873064223040503
$ws.Range("Q24").Value = 75.12425434438133
$ws.Range("R24").Value = 676.118289099432
$ws.Range("S24").Value = 0.001490138532115491
$ws.Range("T24").Value = 0.001490138532115491

$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.481706
$ws.Range("H25").Value = 1.445118
$ws.Range("I25").Value = 0.007955619000060676
$ws.Range("J25").Value = 0.007955619000060676
$ws.Range("M25").Value = 132.7811556666666
$ws.Range("N25").Value = 398.343467
$ws.Range("O25").Value = 0.1594744063806953
$ws.Range("P25").Value = 0.1594744063806954
$ws.Range("Q25").Value = 63.96147937156732
$ws.Range("R25").Value = 575.653314344106
$ws.Range("S25").Value = 0.001268717617425657
$ws.Range("T25").Value = 0.001268717617425658

$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0.481706
$ws.Range("H26").Value = 1.445118
$ws.Range("I26").Value = 0.007955619000060676
$ws.Range("J26").Value = 0.007955619000060676
$ws.Range("M26").Value = 131.9384486666667
$ws.Range("N26").Value = 395.815346
$ws.Range("O26").Value = 0.1584622883741672
$ws.Range("P26").Value = 0.1584622883741672
$ws.Range("Q26").Value = 63.55554235342532
$ws.Range("R26").Value = 571.999881180828
$ws.Range("S26").Value = 0.001260665592182618
$ws.Range("T26").Value = 0.001260665592182618
